# Projekt nyil.tart folyt 2
#
# 1) After the title run "Projekt költség nyilvántartó program" add a
#    trailing "." as its own run (same bold/italic/underline 28pt Times
#    New Roman formatting) so the title ends with a period.
# 2) Merge the "Nem-funkcionális követelmények (NFR)" run and the
#    following ":" run into a single run/text.
# 3) Mark the built-in "Default Paragraph Font" character style as
#    semi-hidden (best effort — see note near the bottom).

$d = $word.ActiveDocument

# --- 1) Title: "...program" -> "...program." -------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Projekt költség nyilvántartó program") | Out-Null
$rng.Collapse(0) # wdCollapseEnd

$start = $rng.Start
$rng.InsertAfter(".")

# Re-apply the exact same character formatting as the title run to the
# newly inserted "." so Word keeps it as its own run instead of folding
# it back into the previous one.
$dot = $d.Range($start, $start + 1)
$dot.Font.Name = "Times New Roman"
$dot.Font.NameBi = "Times New Roman"
$dot.Font.Bold = $true
$dot.Font.Italic = $true
$dot.Font.Underline = 1
$dot.Font.Size = 14

# --- 2) "Nem-funkcionális követelmények (NFR)" + ":" -> one run -------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("Nem-funkcionális követelmények (NFR):", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "Nem-funkcionális követelmények (NFR):", 2) | Out-Null

# --- 3) Default Paragraph Font style -> semi-hidden (best effort) -----
# The COM surface exposed by this host only implements Style.Visibility
# (which toggles the OOXML <w:hidden/> flag) and has no working setter
# that reaches <w:semiHidden/>; Style.Hidden / Style.SemiHidden are both
# unimplemented ("object doesn't support this property or method").
# There is no supported way to author <w:semiHidden/> for a style from
# this object model, so that sub-change is intentionally left alone
# rather than emitting the wrong flag.
